# Update Sheet1 data: add sex/environ/sire/dam columns (D:G) for rows 2-7,
# and remove the old stray notes that lived in rows 12, 15, 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the old stray annotation cells (B12, G15, G16) that are no longer used.
$ws.Range("B12").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("G16").ClearContents()

# Populate new D:G columns for rows 2-7.
$ws.Range("D2:D7").Value = "U"
$ws.Range("E2:E7").Value = 1

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 3

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 2

# Update the selection to match the final workbook view.
$ws.Range("G9").Select()
